$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 3 (the second data row) - shifts everything below it up by one.
$ws.Rows.Item(3).Delete()
